$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '66.161.18'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  +2.35%  '

$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.235.52'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  +5.85%  '

$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  -0.10%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '579.52'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  +4.33%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '151.93'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  +6.99%  '

$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  -0.01%  '

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '3.227.12'
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  +5.80%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.514'
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  +4.69%  '

$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  +8.50%  '

$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  +5.44%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.489'
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  +5.40%  '

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '38.35'
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  +5.96%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.0000234'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  +5.69%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '3.757.49'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  +6.09%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '66.186.16'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  +2.23%  '

$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '543.28'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  +11.45%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '3.233.38'

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.115'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  +2.91%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '7.14'
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  +7.11%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '14.62'
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  +7.01%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.744'
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  +8.20%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '7.81'
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  +8.99%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '13.57'
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  +7.46%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '81.36'
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  +3.02%  '

$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  -0.13%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '9.37'
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  +20.21%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.96'
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  +8.05%  '

$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  +7.64%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '27.69'
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  +6.63%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '2.76'
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  +5.58%  '

$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  +0.14%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.18'
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  +5.24%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '564.72'
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  +2.57%  '

$ws.Range('B35').NumberFormat = '@'
$ws.Range('B35').Value = 'NEARProtocol'
$ws.Range('C35').NumberFormat = '@'
$ws.Range('C35').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '5.65'
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  +3.91%  '

$ws.Range('B36').NumberFormat = '@'
$ws.Range('B36').Value = 'Filecoin'
$ws.Range('C36').NumberFormat = '@'
$ws.Range('C36').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '6.37'
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  +7.18%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.0459'
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  +9.82%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '54.65'
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  +3.37%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.0863'
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  +7.32%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.130'
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  +7.07%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '3.205.51'
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  +10.66%  '

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '2.91'
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  +3.64%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '8.61'
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  +4.31%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.285'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  +16.33%  '

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '2.35'
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  +12.52%  '

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '26.55'
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  +6.99%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '1.00'
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  +0.05%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.0₃0557'
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  +3.56%  '

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '124.58'
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  +3.15%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.113'
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  +3.80%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '2.21'
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  +8.07%  '
